# Update "想去人数" (wish-to-attend count) figures in column F on the
# "展览" (rId1 / sheet1) and "全部类型" (rId4 / sheet4) worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 301
$wsExhibit.Range("F4").Value = 8038
$wsExhibit.Range("F5").Value = 5861
$wsExhibit.Range("F8").Value = 14
$wsExhibit.Range("F10").Value = 287
$wsExhibit.Range("F11").Value = 390

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 301
$wsAll.Range("F4").Value = 8038
$wsAll.Range("F5").Value = 5861
$wsAll.Range("F8").Value = 14
$wsAll.Range("F10").Value = 287
$wsAll.Range("F14").Value = 390
